# Fix "Recorded By" (column G) ordering on the "Session Analysis Results" sheet.
# The author names/emails in each comma-separated cell were reordered
# (reversed) so that the human-looking name/backup address comes first
# instead of "System"/"system" leading the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    2   = "backup@backdoor.com, System, system"
    4   = "backup@backdoor.com, System"
    5   = "backup@backdoor.com, System"
    8   = "backup@backdoor.com, System"
    11  = "dnasr281@gmail.com, System"
    17  = "dnasr281@gmail.com, System"
    29  = "backup@backdoor.com, System, system"
    31  = "backup@backdoor.com, System"
    32  = "backup@backdoor.com, System"
    35  = "backup@backdoor.com, System"
    38  = "dnasr281@gmail.com, System"
    44  = "dnasr281@gmail.com, System"
    56  = "backup@backdoor.com, System, system"
    58  = "backup@backdoor.com, System"
    59  = "backup@backdoor.com, System"
    62  = "backup@backdoor.com, System"
    65  = "dnasr281@gmail.com, System"
    71  = "dnasr281@gmail.com, System"
    83  = "backup@backdoor.com, System"
    84  = "backup@backdoor.com, System"
    85  = "backup@backdoor.com, System"
    90  = "dnasr281@gmail.com, admin@admin.com"
    96  = "dnasr281@gmail.com, System"
    97  = "dnasr281@gmail.com, System"
    99  = "dnasr281@gmail.com, System"
    109 = "backup@backdoor.com, System"
    110 = "backup@backdoor.com, System"
    111 = "backup@backdoor.com, System"
    116 = "dnasr281@gmail.com, admin@admin.com"
    122 = "dnasr281@gmail.com, System"
    123 = "dnasr281@gmail.com, System"
    125 = "dnasr281@gmail.com, System"
    135 = "backup@backdoor.com, System"
    136 = "backup@backdoor.com, System"
    137 = "backup@backdoor.com, System"
    142 = "dnasr281@gmail.com, admin@admin.com"
    148 = "dnasr281@gmail.com, System"
    149 = "dnasr281@gmail.com, System"
    151 = "dnasr281@gmail.com, System"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
